$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (D29): PAIR troubleshooting steps - wording tweaks (조사하기/분석하기/분리하기/해결하기)
$ws.Range("D29").Value = '<a style="color: red; font-weight: bold;">Volvo Internal Use Only</a>
<br>
<br>이 문서는 문제를 해결하기 위한 기본 단계에 대해 설명합니다.
<br>
<br>고객의 문제를 해결하기 위해 PAIR를 적용해야 합니다. PAIR를 통해 정확한 원인을 확인하고 분석하여 정확한 원인을 도출할 수 있습니다.
<br>
<br><strong>P: Probe (조사하기)</strong>
<br>
<br>문제를 조사하고 분석하는 단계입니다. 문제의 본질을 이해하기 위해 다양한 질문을 던지고, 문제의 범위와 영향을 받는 사람들을 파악합니다. 이 단계에서는 문제의 초기 원인과 상황을 파악하기 위해 데이터를 수집하고, 관찰하는 단계를 수행합니다.
<br>
<br><strong>A: Analyze (분석하기)</strong>
<br>
<br>수집된 데이터를 바탕으로 문제의 원인과 구조를 분석하는 단계입니다. 문제의 근본 원인을 파악하고, 문제 해결을 위해 필요한 정보를 체계적으로 분석합니다. 이 과정에서 다양한 분석 도구와 기법을 사용하여 데이터를 해석하고, 문제의 핵심 요소를 도출합니다.
<br>
<br><strong>I: Isolate (분리하기)</strong>
<br>
<br>문제의 원인을 분리하고 격리하는 단계입니다. 문제를 더 작은 부분으로 나누어 각 부분을 독립적으로 분석하고, 특정 원인을 식별하여 격리합니다. 이를 통해 문제를 더 구체적이고 세부적으로 다루며, 각 원인별로 맞춤형 해결책을 모색합니다.
<br>
<br><strong>R: Resolve (해결하기)</strong>
<br>
<br>문제를 해결하고, 해결책을 실행하는 단계입니다. 분석과 격리 단계에서 도출된 해결책을 실행하여 문제를 해결합니다. 또한, 해결책의 효과를 검토하고, 지속적인 모니터링을 통해 문제 재발 방지 및 추가적인 개선을 도모합니다.'

# Row 32 (D32): Volvo Cars app aircon issue - fix unclosed <strong> tag typo
$ws.Range("D32").Value = '고객이 Volvo Cars앱에서 예기치 않게 공조 장치를 조정할 수 없다고 보고할 수 있습니다.
<br>
<br><strong>차량에 로그인된 Volvo ID와 휴대폰에 설치된 Volvo Cars앱의 Volvo ID가 서로 다른 경우</strong>
<br>1. 차량의 우측 하단의 설정으로 이동합니다.
<br>2. 프로필로 이동하여 프로필 정보에 있는 계정을 누릅니다.
<br>3. Volvo Cars앱과 동일한 계정으로 로그인되어있는지 확인합니다. 다른 경우 삭제 후 로그인을 진행합니다.
<br>4. 그 다음 설정에서 프로필로 이동합니다. 이후 Volvo Cars 앱 장치 항목으로 이동합니다.
<br>5. 전화 연결 항목을 선택하고 안내를 따릅니다.
<br>
<br><strong>동일한 계정이 로그인된 경우</strong>
<br>1. 차량 상단바를 밑으로 내려 알림 센터를 엽니다.
<br>2. 차량 프로필을 오너로 변경해봅니다.'

# Row 44: replace "마케팅 문자 수신 거부 요청" with the expanded opt-out + info removal howto
$ws.Range("A44").Value = '마케팅 수신 거부 요청 및 정보 제거 처리하기'
$ws.Range("B44").Value = 45461
$ws.Range("D44").Value = '<br>고객이 마케팅 수신 거부를 요청하거나, 정보 제거를 원하는 경우 아래 단계에 따라 진행합니다.
<br>
<br><h3>마케팅 수신 거부를 요청하는 경우</h3>
<br><strong>1.</strong> DMS에서 고객 전화번호+이름을 조회 후 수정을 누릅니다.
<br><strong>2.</strong> 고객 마케팅 동의여부 수정에서 거부로 변경합니다.
<br>
<br><h3>정보 제거를 원하는 경우</h3>
<br><strong>1.</strong> 고객명과 전화번호를 수집합니다.
<br><strong>2.</strong> 아래 메일 주소로 수신 거부 요청을 발송합니다.
<br>
<br><strong>TO:</strong> heejaekelly.byun@volvocars.com
<br><strong>CC:</strong> yhj@mytable.co.kr / volvo@mytable.co.kr / jungsoo.kim.2@volvocars.com / eunha.jung@volvocars.com / dohyeon.kim@volvocars.com / seungjin.lee@volvocars.com / jeongho.suh@volvocars.com / heemin.park@volvocars.com / hyunho.shim@volvocars.com'

# Row height for row 44 grew to fit the longer content
$ws.Rows.Item(44).RowHeight = 224.25

# Update the saved view state to match where the edit left the selection
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C44").Select()
